$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet grows a new "pt_max" column (value 50 for every data row),
# inserted right before the existing "boson" column. Everything from the
# old column E onward shifts one column to the right (E->F, F->G, ... M->N).
#
# Rename the about-to-shift "syst_c" header to "syst_u" while it is still
# sitting in its original column (H). Doing the rename before the insert
# means the new shared-string entry lands in the same relative slot
# (right after "norm_c") as in the authored workbook.
$ws.Range("H1").Value = "syst_u"

# Insert the new column before the old column E ("boson"); this shifts
# E:M -> F:N, updates column widths/formatting, the used range, etc.,
# exactly like using Excel's "Insert Column" on column E.
$ws.Range("E1").EntireColumn.Insert()

# Fill in the freshly inserted column E: a header label plus a constant
# value of 50 for each of the 16 data rows.
$ws.Range("E1").Value = "pt_max"
$ws.Range("E2:E17").Value = 50

# Match the saved selection/active cell.
$ws.Range("E2:E17").Select()

# Best-effort: reposition the window to match the authored file's saved
# window coordinates (not guaranteed to round-trip through every host,
# but harmless to attempt).
$win = $excel.ActiveWindow
$win.Left = 0
$win.Top = 500
